# Daily attendance processing - 2025-11-14 15:46:30
# Swap the first and last entries of the comma-separated "Recorded By"
# list (column G) for every data row that has more than one entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row  # xlUp
if ($lastRow -lt 2) { $lastRow = 2 }

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($null -eq $value) { continue }

    $text = [string]$value
    if ($text -notmatch ',') { continue }

    $parts = $text -split ',\s*'
    $n = $parts.Count
    if ($n -lt 2) { continue }

    $tmp = $parts[0]
    $parts[0] = $parts[$n - 1]
    $parts[$n - 1] = $tmp

    $cell.Value2 = [string]::Join(', ', $parts)
}
